$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("keys")

# Move the "<missing>" markers from column E into column D for rows 4, 7, 9
# (D was empty there before), preserving the highlighted "missing" style.
foreach ($r in 4, 7, 9) {
    $dst = $ws.Cells.Item($r, 4)   # D
    $src = $ws.Cells.Item($r, 5)   # E
    $dst.Value = "<missing>"
    $src.Copy()
    $dst.PasteSpecial(-4122)       # xlPasteFormats
    $dst.Value = "<missing>"
    $src.ClearContents()
}

# Add the new "<missing>" marker cells in column B for rows 5, 8, 10, matching
# the format already used by the other "<missing>" cells in those rows.
foreach ($r in 5, 8, 10) {
    $fmt = $ws.Cells.Item($r, 1)   # A (already styled "<missing>" cell)
    $dst = $ws.Cells.Item($r, 2)   # B
    $fmt.Copy()
    $dst.PasteSpecial(-4122)       # xlPasteFormats
    $dst.Value = "<missing>"
}

# Make the C5/C8/C10 number-cell formatting match the rest of the
# "<missing>"-row formatting (same visual style; de-duplicates the xf).
foreach ($r in 5, 8, 10) {
    $fmt = $ws.Cells.Item($r, 4)   # D (already on the shared style)
    $dst = $ws.Cells.Item($r, 3)   # C
    $val = $dst.Value2
    $fmt.Copy()
    $dst.PasteSpecial(-4122)       # xlPasteFormats
    $dst.Value = $val
}

# Remove the now-empty, now-unused column E entirely (shifts nothing left of
# it; there is nothing to its right).
$ws.Columns.Item(5).Delete()

# Shrink the AutoFilter range to A1:D11.
$ws.AutoFilterMode = $false
$ws.Range("A1:D11").AutoFilter() | Out-Null

# Update the _FilterDatabase defined name to the new range.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=keys!`$A`$1:`$D`$11"
